$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) stays the same text, just rewritten for clarity ----
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "created_at"
$ws.Cells.Item(1,3).Value = "first_name"
$ws.Cells.Item(1,4).Value = "last_name"
$ws.Cells.Item(1,5).Value = "email"
$ws.Cells.Item(1,6).Value = "password"
$ws.Cells.Item(1,7).Value = "recovery_question"
$ws.Cells.Item(1,8).Value = "recovery_answer"

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = "53af4926-52ee-41d0-9acc-ae7230000001"
$ws.Cells.Item(2,2).Value = "2020-03-25 02:17:06"
$ws.Cells.Item(2,3).Value = "Aina"
$ws.Cells.Item(2,4).Value = "Jesulayomi"
$ws.Cells.Item(2,5).Value = "jesulayomi@schub.com"
$ws.Cells.Item(2,6).Value = "ajpwd"
$ws.Cells.Item(2,7).Value = "What is your baby's name"
$ws.Cells.Item(2,8).Value = "Micoliser"

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = "53af4926-52ee-41d0-9acc-ae7230000002"
$ws.Cells.Item(3,2).Value = "2020-03-25 02:17:06"
$ws.Cells.Item(3,3).Value = "Iwelumo"
$ws.Cells.Item(3,4).Value = "Samuel"
$ws.Cells.Item(3,5).Value = "samuel@schub.com"
$ws.Cells.Item(3,6).Value = "ispwd"
$ws.Cells.Item(3,7).Value = "What is your baby's name"
$ws.Cells.Item(3,8).Value = "Jesulayomi"

# ---- Row 4 (Beta Tester -> Bot / Beta-Tester) ----
$ws.Cells.Item(4,1).Value = "53af4926-52ee-41d0-9acc-ae7230000003"
$ws.Cells.Item(4,2).Value = "2023-03-25 02:17:07"
$ws.Cells.Item(4,3).Value = "Beta-Tester"
$ws.Cells.Item(4,4).Value = "Bot"
$ws.Cells.Item(4,5).Value = "tester@schub.com"
$ws.Cells.Item(4,6).Value = "testerpwd"
$ws.Cells.Item(4,7).Value = "Hello there?"
$ws.Cells.Item(4,8).Value = "General Kenobi"

# ---- Row 5 (new - Alfred Tuva) ----
$ws.Cells.Item(5,1).Value = "53af4926-52ee-41d0-9acc-ae7230000004"
$ws.Cells.Item(5,2).Value = "2023-03-25 02:17:08"
$ws.Cells.Item(5,3).Value = "Alfred"
$ws.Cells.Item(5,4).Value = "Tuva"
$ws.Cells.Item(5,5).Value = "alfred@schub.com"
$ws.Cells.Item(5,6).Value = "alfredpwd"
$ws.Cells.Item(5,7).Value = "Hello there?"
$ws.Cells.Item(5,8).Value = "General Kenobi"

# ---- Row 6 (new - Martins Ndifon) ----
$ws.Cells.Item(6,1).Value = "53af4926-52ee-41d0-9acc-ae7230000005"
$ws.Cells.Item(6,2).Value = "2023-03-25 02:17:09"
$ws.Cells.Item(6,3).Value = "Martins"
$ws.Cells.Item(6,4).Value = "Ndifon"
$ws.Cells.Item(6,5).Value = "martins@schub.com"
$ws.Cells.Item(6,6).Value = "martinspwd"
$ws.Cells.Item(6,7).Value = "Hello there?"
$ws.Cells.Item(6,8).Value = "General Kenobi"

# ---- Row 7 (new - Julien Barbier) ----
$ws.Cells.Item(7,1).Value = "53af4926-52ee-41d0-9acc-ae7230000006"
$ws.Cells.Item(7,2).Value = "2023-03-25 02:17:10"
$ws.Cells.Item(7,3).Value = "Julien"
$ws.Cells.Item(7,4).Value = "Barbier"
$ws.Cells.Item(7,5).Value = "julien@schub.com"
$ws.Cells.Item(7,6).Value = "julienpwd"
$ws.Cells.Item(7,7).Value = "Hello there?"
$ws.Cells.Item(7,8).Value = "General Kenobi"

# ---- Hyperlinks for the new email cells ----
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:alfred@schub.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:martins@schub.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:julien@schub.com")

# ---- Selection / view ----
$ws.Range("F10").Select()
